$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, reflecting repulled / recalculated data.
$values = @{
    2  = -2
    4  = -4
    5  = -3
    6  = -3
    7  = -4
    10 = 1
    12 = -5
    13 = -3
    14 = 2
    15 = 2
    16 = -1
    17 = 2
    18 = 3
    19 = -3
    20 = 1
    21 = -6
    22 = 4
    23 = -4
    24 = 11
    25 = 1
    26 = 4
    27 = -3
    29 = 2
    30 = 4
    31 = 4
    32 = 6
    33 = -4
    34 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
